# add notification for wrong raw occurrence dataset
#
# The raw-occurrence termCode recorded for the "Potatoes {Mixing}" process
# facet was wrong (missing a distinguishing suffix), and the derived
# LB/MB/UB columns for the three potato process-facet rows (Roasting,
# Frying, Mixing) were using the wrong scaling factor. Fix both.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the mislabeled raw occurrence termCode for "Potatoes {Mixing}" ---
$ws.Range("A23").Value = "A00ZT#F28.A0CRlss"

# --- Correct the derived-value scaling factor (0.7x -> 2x) ---
# Set the anchor (E) column cell by cell and the trailing F:J columns as a
# row range, matching the original per-row formula layout.
$ws.Range("E21").Formula = "=E12*2"
$ws.Range("F21:J21").Formula = "=F12*2"

$ws.Range("E22").Formula = "=E12*2"
$ws.Range("F22:J22").Formula = "=F12*2"

$ws.Range("E23").Formula = "=E13*2"
$ws.Range("F23:J23").Formula = "=F13*2"

# --- Move the selection to the range that was just edited ---
$ws.Range("E23:J23").Select()
